$d = $word.ActiveDocument

# --- Paragraph 1 ("Journal de bord :") ---------------------------------
# The visible text does not change, but the stray proofing markup
# (proofErr spellStart/gramStart/spellEnd) around "bord" is cleaned up
# and the "nal de " / "bord" runs merge into a single run. Replacing the
# matched text with itself achieves this run-merge/tag-cleanup without
# altering the paragraph's actual wording.
$rng1 = $d.Paragraphs(1).Range
$rng1.Find.ClearFormatting()
$rng1.Find.Execute("nal de bord", $true, $false, $false, $false, $false, $true, 1, $false, "nal de bord", 2) | Out-Null

# --- Paragraph 2 (journal entry) ----------------------------------------
# Append the new journal text. The existing "_GoBack" bookmark is removed
# from its old position and re-created between "Test s" and
# "ur bt-200, debut Gantt." in the newly appended text.
$p2 = $d.Paragraphs(2)

$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

$insertPoint = $p2.Range
$insertPoint.Collapse(0)
$insertPoint.InsertAfter("Test sur bt-200, debut Gantt.")

$splitPos = $p2.Range.End - 1 - "ur bt-200, debut Gantt.".Length
$bookmarkRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
